# Added population and income to route dataset
# Adds a "Locid" (airport/location code) column E to the Categorization sheet.

$wb = $excel.ActiveWorkbook
$income = $wb.Worksheets.Item("Income")
$cat = $wb.Worksheets.Item("Categorization")

# --- Header (E1) ---
$cat.Range("E1").Value = "Locid"

# --- Data rows (E2:E31): location codes matching each route's airport ---
$codes = @{
    2  = "ATL"
    3  = "LAX"
    4  = "ORD"
    5  = "DFW"
    6  = "JFK"
    7  = "DEN"
    8  = "SFO"
    9  = "CLT"
    10 = "LAS"
    11 = "PHX"
    12 = "IAH"
    13 = "MIA"
    14 = "SEA"
    15 = "EWR"
    16 = "MCO"
    17 = "MSP"
    18 = "DTW"
    19 = "BOS"
    20 = "PHL"
    21 = "LGA"
    22 = "FLL"
    23 = "BWI"
    24 = "IAD"
    25 = "MDW"
    26 = "SLC"
    27 = "DCA"
    28 = "HNL"
    29 = "SAN"
    30 = "TPA"
    31 = "PDX"
}

# Rows that use the "highlighted" (red) label style, matching the Income sheet's
# equivalent rows (e.g. Fort Lauderdale / Honolulu).
$highlightRows = @(22, 28)

# Pick up the two pre-existing cell formats from the Income sheet (column A)
# so the new column matches the rest of the workbook's look & feel.
$income.Range("A2").Copy()
foreach ($r in 2..31) {
    if ($highlightRows -contains $r) { continue }
    $cat.Range("E$r").PasteSpecial(-4122)
}

$income.Range("A23").Copy()
foreach ($r in $highlightRows) {
    $cat.Range("E$r").PasteSpecial(-4122)
}

foreach ($r in 2..31) {
    $cat.Range("E$r").Value = $codes[$r]
}

# --- View-state tweaks ---
$cat.Range("E7").Select()
$income.Range("A10").Select()
